$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the "zh-cn"/"de-de" status columns (17.216 -> 13.410 chars) ---
$overview.Range("E:F").ColumnWidth = 12.5
$zhcn.Range("C:C").ColumnWidth = 12.5
$dede.Range("C:C").ColumnWidth = 12.5
